$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 73.8161485502559
$ws.Range("C2").Value = 72.4542861110033
$ws.Range("D2").Value = 75.1780109895084
$ws.Range("B3").Value = 43.8285876250807
$ws.Range("C3").Value = 36.981532036591
$ws.Range("D3").Value = 50.6756432135703
$ws.Range("C4").Value = 70.8472369591563
$ws.Range("D4").Value = 80.062499290752
$ws.Range("C5").Value = 89.1248437720431
$ws.Range("D5").Value = 94.1999973806463
$ws.Range("C6").Value = 69.2225897033517
$ws.Range("D6").Value = 77.6714823287634
$ws.Range("B7").Value = 64.8131869221429
$ws.Range("C7").Value = 59.3428029309974
$ws.Range("D7").Value = 70.2835709132883
$ws.Range("C8").Value = 61.4234131747877
$ws.Range("D8").Value = 76.081554069256
$ws.Range("C9").Value = 72.7430653096427
$ws.Range("D9").Value = 86.9947356115996
$ws.Range("C10").Value = 65.9680421577375
$ws.Range("D10").Value = 77.5928469114932
$ws.Range("C11").Value = 77.227517680188
$ws.Range("D11").Value = 90.9665315268747
$ws.Range("C12").Value = 75.6441025705525
$ws.Range("D12").Value = 94.3353077679108
$ws.Range("C13").Value = 68.7502708550797
$ws.Range("D13").Value = 76.0748827650365
$ws.Range("C14").Value = 76.5690856435808
$ws.Range("D14").Value = 83.2423600984255
$ws.Range("C15").Value = 46.1259689996369
$ws.Range("D15").Value = 55.1860069092131
$ws.Range("B16").Value = 69.211130206042
$ws.Range("C16").Value = 67.7830711918507
$ws.Range("D16").Value = 70.6391892202333
$ws.Range("B17").Value = 43.9777455537456
$ws.Range("C17").Value = 37.0123347353316
$ws.Range("D17").Value = 50.9431563721596
$ws.Range("C18").Value = 66.1454151446879
$ws.Range("D18").Value = 75.9383056859227
$ws.Range("C19").Value = 86.5789684031514
$ws.Range("D19").Value = 92.2234242404618
$ws.Range("C20").Value = 60.8006265400238
$ws.Range("D20").Value = 70.199568787608
$ws.Range("C21").Value = 47.8901285995169
$ws.Range("D21").Value = 58.8550560842209
$ws.Range("C22").Value = 50.2806013246459
$ws.Range("D22").Value = 65.7189964020433
$ws.Range("C23").Value = 64.8845068373116
$ws.Range("D23").Value = 80.7175234088101
$ws.Range("C24").Value = 63.3460230300396
$ws.Range("D24").Value = 75.4339258938102
$ws.Range("C25").Value = 76.0339981863004
$ws.Range("D25").Value = 89.8853810009777
$ws.Range("C26").Value = 74.3329265016726
$ws.Range("D26").Value = 93.2948213208992
$ws.Range("B27").Value = 67.612743921124
$ws.Range("C27").Value = 63.8143400774096
$ws.Range("D27").Value = 71.4111477648383
$ws.Range("C28").Value = 76.6430556907113
$ws.Range("D28").Value = 83.2787814520089
$ws.Range("C29").Value = 45.4719434410811
$ws.Range("D29").Value = 54.6334088578313
$ws.Range("B30").Value = 68.2106860437125
$ws.Range("C30").Value = 66.767714507923
$ws.Range("D30").Value = 69.6536575795019
$ws.Range("B31").Value = 40.0825994382391
$ws.Range("C31").Value = 33.0521475410672
$ws.Range("D31").Value = 47.1130513354109
$ws.Range("C32").Value = 65.5589503981622
$ws.Range("D32").Value = 75.4602996898891
$ws.Range("C33").Value = 82.6063958744822
$ws.Range("D33").Value = 89.0125191611756
$ws.Range("C34").Value = 60.7076351601826
$ws.Range("D34").Value = 70.1994269609814
$ws.Range("C35").Value = 53.5359954959464
$ws.Range("D35").Value = 64.5634501788632
$ws.Range("C36").Value = 44.5445750246031
$ws.Range("D36").Value = 61.3949180286547
$ws.Range("C37").Value = 59.7618877956136
$ws.Range("D37").Value = 76.4648893365579
$ws.Range("C38").Value = 56.3357986578306
$ws.Range("D38").Value = 68.6509199418123
$ws.Range("C39").Value = 75.4596634190305
$ws.Range("D39").Value = 88.9118906290386
$ws.Range("B40").Value = 80.9588089320727
$ws.Range("C40").Value = 71.4176457256646
$ws.Range("D40").Value = 90.4999721384807
$ws.Range("C41").Value = 61.957440302649
$ws.Range("D41").Value = 69.5663846767349
$ws.Range("C42").Value = 74.6238068681753
$ws.Range("D42").Value = 81.5525407490971
$ws.Range("C43").Value = 44.1396772251881
$ws.Range("D43").Value = 52.7926027872294
